# Fix: append the missing end-of-log rows that point Raul's log back at
# real data (the bug had the event pointing at a null/blank reference).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

# r, A(Staff), B(Date serial), C(Time), D(Building), E(Room), F(Notes), wrapped(bool)
$rows = @(
    @(770, "Pickup Mic",       42666, "1630", "OSG", "1005", "Return one podium,  2 desk and 2 audience mics with cables and stands to booth behind stage", $true),
    @(771, "AV Shutdown",      42666, "1630", "OSG", "1005", $null, $false),
    @(772, "AV Shutdown",      42666, "1700", "SSB", "E111", $null, $false),
    @(773, "AV Shutdown",      42666, "2100", "SSB", "W141", $null, $false),
    @(774, "Pickup Mic",       42666, "1630", "OSG", "2001", "Return 2 desk mics, mixer, speaker , stands and cables to OSG 1014L", $true),
    @(775, "Pickup Mic",       42666, "1630", "OSG", "2002", "Return 2 desk mics, mixer, stands and cables to OSG 1014L", $false),
    @(776, "Pickup Mic",       42666, "1630", "OSG", "1014", "Return 2 desk and one audience mic , stands , cables and mixer to OSG 1014L", $true),
    @(777, "Pickup Skype Kit", 42666, "1630", "OSG", "1014", "Return with USB extenders to OSG 1014L", $false),
    @(778, "AV Shutdown",      42666, "1630", "OSG", "1014", $null, $false)
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    $ws.Range("A$rowNum").Value = $r[1]
    $ws.Range("B$rowNum").Value = $r[2]
    $ws.Range("C$rowNum").Value = $r[3]
    $ws.Range("D$rowNum").Value = $r[4]
    $ws.Range("E$rowNum").Value = $r[5]
    if ($r[6]) {
        $ws.Range("F$rowNum").Value = $r[6]
    }
    if ($r[7]) {
        $ws.Rows($rowNum).RowHeight = 30
    }
}

$ws.Range("A778").Select()
